# Excel COM-interop script: add syst0_c/syst1_c/syst2_c columns (M:O) and
# rename the existing syst_u column (H) to syst_tot.
#
# The underlying per-row numeric "syst_u"/"syst_tot" values in column H are
# unchanged; only the header label moves from "syst_u" to "syst_tot".
# Three brand new columns are appended (M, N, O) holding the individual
# systematic-uncertainty components (syst0_c, syst1_c, syst2_c) whose
# quadrature sum reproduces column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header H1: "syst_u" -> "syst_tot" (values in H2:H7 untouched) ---
$ws.Range("H1").Value = "syst_tot"

# --- New header cells M1:O1, styled centered/top-aligned, 11pt Calibri ---
$headerRange = $ws.Range("M1:O1")
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11

$ws.Range("M1").Value = "syst0_c"
$ws.Range("N1").Value = "syst1_c"
$ws.Range("O1").Value = "syst2_c"

# --- New numeric data for rows 2-7 ---
$ws.Range("M2").Value = 0.0059
$ws.Range("N2").Value = 0.005
$ws.Range("O2").Value = 0.0023

$ws.Range("M3").Value = 0.0039
$ws.Range("N3").Value = 0.0048
$ws.Range("O3").Value = 0.0029

$ws.Range("M4").Value = 0.0092
$ws.Range("N4").Value = 0.005
$ws.Range("O4").Value = 0.0034

$ws.Range("M5").Value = 0.0072
$ws.Range("N5").Value = 0.0048
$ws.Range("O5").Value = 0.004

$ws.Range("M6").Value = 0.0081
$ws.Range("N6").Value = 0.005
$ws.Range("O6").Value = 0.0053

$ws.Range("M7").Value = 0.0117
$ws.Range("N7").Value = 0.0042
$ws.Range("O7").Value = 0.0058

# --- Column width / selection to match the refreshed layout ---
$ws.Columns("M").ColumnWidth = 15.1640625
$ws.Range("L17").Select() | Out-Null
